$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("展览")
$ws.Range("F2").Value = 203
$ws.Range("F6").Value = 269
$ws.Range("F12").Value = 114
$ws.Range("F13").Value = 2422
$ws.Range("F14").Value = 39
$ws.Range("F16").Value = 11
$ws.Range("F18").Value = 530
$ws.Range("F19").Value = 569
$ws.Range("F22").Value = 49
$ws.Range("F24").Value = 1987
$ws.Range("F25").Value = 4115
$ws.Range("F27").Value = 63
$ws.Range("F28").Value = 1199
$ws.Range("F30").Value = 2108
$ws.Range("F32").Value = 471
$ws.Range("F34").Value = 123
$ws.Range("F36").Value = 428
$ws.Range("F41").Value = 423

$ws = $wb.Worksheets.Item("演出")
$ws.Range("F2").Value = 39

$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F2").Value = 203
$ws.Range("F6").Value = 269
$ws.Range("F12").Value = 114
$ws.Range("F13").Value = 2422
$ws.Range("F14").Value = 39
$ws.Range("F16").Value = 39
$ws.Range("F17").Value = 11
$ws.Range("F19").Value = 530
$ws.Range("F20").Value = 569
$ws.Range("F23").Value = 49
$ws.Range("F25").Value = 1987
$ws.Range("F26").Value = 4115
$ws.Range("F28").Value = 63
$ws.Range("F29").Value = 1199
$ws.Range("F31").Value = 2108
$ws.Range("F33").Value = 471
$ws.Range("F35").Value = 123
$ws.Range("F37").Value = 428
$ws.Range("F42").Value = 423

